$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1846153846153846
$ws.Range("C2").Value = 0.5435897435897435
$ws.Range("J2").Value = 0.01025641025641026
$ws.Range("P2").Value = 0.1641025641025641
$ws.Range("S2").Value = 0.09743589743589744
$ws.Range("C3").Value = 0.01851851851851852
$ws.Range("J3").Value = 0.04629629629629629
$ws.Range("P3").Value = 0.7314814814814815
$ws.Range("S3").Value = 0.2037037037037037
$ws.Range("J4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.6744186046511628
$ws.Range("S4").Value = 0.3023255813953488
$ws.Range("B6").Value = 0.05
$ws.Range("D6").Value = 0.015
$ws.Range("F6").Value = 0.04
$ws.Range("J6").Value = 0.22
$ws.Range("O6").Value = 0.005
$ws.Range("Q6").Value = 0.11
$ws.Range("R6").Value = 0.065
$ws.Range("S6").Value = 0.495
$ws.Range("B7").Value = 0.05294117647058823
$ws.Range("D7").Value = 0.04117647058823529
$ws.Range("E7").Value = 0.005882352941176471
$ws.Range("F7").Value = 0.05882352941176471
$ws.Range("J7").Value = 0.1176470588235294
$ws.Range("O7").Value = 0.01764705882352941
$ws.Range("Q7").Value = 0.1705882352941177
$ws.Range("R7").Value = 0.1352941176470588
$ws.Range("S7").Value = 0.4
$ws.Range("B8").Value = 0.062
$ws.Range("D8").Value = 0.014
$ws.Range("F8").Value = 0.078
$ws.Range("J8").Value = 0.08599999999999999
$ws.Range("O8").Value = 0.026
$ws.Range("Q8").Value = 0.174
$ws.Range("R8").Value = 0.092
$ws.Range("S8").Value = 0.468
$ws.Range("B9").Value = 0.06349206349206349
$ws.Range("D9").Value = 0.005291005291005291
$ws.Range("F9").Value = 0.0582010582010582
$ws.Range("J9").Value = 0.08465608465608465
$ws.Range("O9").Value = 0.01587301587301587
$ws.Range("Q9").Value = 0.1428571428571428
$ws.Range("R9").Value = 0.08994708994708994
$ws.Range("S9").Value = 0.5396825396825397
$ws.Range("B10").Value = 0.07860262008733625
$ws.Range("D10").Value = 0.02358078602620087
$ws.Range("F10").Value = 0.06200873362445415
$ws.Range("J10").Value = 0.1240174672489083
$ws.Range("O10").Value = 0.01135371179039301
$ws.Range("Q10").Value = 0.2034934497816594
$ws.Range("R10").Value = 0.07860262008733625
$ws.Range("S10").Value = 0.4183406113537118
$ws.Range("G11").Value = 0.1604938271604938
$ws.Range("J11").Value = 0.06172839506172839
$ws.Range("K11").Value = 0.1934156378600823
$ws.Range("L11").Value = 0.5720164609053497
$ws.Range("S11").Value = 0.01234567901234568
$ws.Range("G12").Value = 0.7551020408163265
$ws.Range("J12").Value = 0.163265306122449
$ws.Range("K12").Value = 0.01360544217687075
$ws.Range("L12").Value = 0.02040816326530612
$ws.Range("S12").Value = 0.04761904761904762
$ws.Range("G13").Value = 0.7222222222222222
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.02777777777777778
$ws.Range("F15").Value = 0.02325581395348837
$ws.Range("H15").Value = 0.1686046511627907
$ws.Range("I15").Value = 0.1279069767441861
$ws.Range("J15").Value = 0.3081395348837209
$ws.Range("K15").Value = 0.04651162790697674
$ws.Range("M15").Value = 0.01744186046511628
$ws.Range("N15").Value = 0.005813953488372093
$ws.Range("O15").Value = 0.02906976744186046
$ws.Range("S15").Value = 0.2732558139534884
$ws.Range("F16").Value = 0.007407407407407408
$ws.Range("H16").Value = 0.237037037037037
$ws.Range("I16").Value = 0.1407407407407407
$ws.Range("J16").Value = 0.3407407407407407
$ws.Range("K16").Value = 0.1037037037037037
$ws.Range("N16").Value = 0.007407407407407408
$ws.Range("O16").Value = 0.05185185185185185
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.02827763496143959
$ws.Range("H17").Value = 0.2287917737789203
$ws.Range("I17").Value = 0.115681233933162
$ws.Range("J17").Value = 0.3933161953727506
$ws.Range("K17").Value = 0.08483290488431877
$ws.Range("M17").Value = 0.01542416452442159
$ws.Range("O17").Value = 0.05141388174807198
$ws.Range("S17").Value = 0.08226221079691516
$ws.Range("F18").Value = 0.01075268817204301
$ws.Range("H18").Value = 0.1881720430107527
$ws.Range("I18").Value = 0.08602150537634409
$ws.Range("J18").Value = 0.4193548387096774
$ws.Range("K18").Value = 0.08064516129032258
$ws.Range("M18").Value = 0.01612903225806452
$ws.Range("O18").Value = 0.04838709677419355
$ws.Range("S18").Value = 0.1505376344086022
$ws.Range("F19").Value = 0.01066260472201066
$ws.Range("H19").Value = 0.2444782939832445
$ws.Range("I19").Value = 0.06549885757806551
$ws.Range("J19").Value = 0.38994668697639
$ws.Range("K19").Value = 0.09367859862909368
$ws.Range("M19").Value = 0.01827875095201828
$ws.Range("O19").Value = 0.05788271134805788
$ws.Range("S19").Value = 0.1195734958111196
